# 20210218 管理単位にBI を追加 packing.py
# Updates the shipment-results sheet: rows 5 & 6 get new shipment data
# (2021-02-19 / 高圧化工 / UV-PP items sent to Osaka), row 7 (the old
# タキ倉庫 line) is removed entirely, and a handful of columns are
# resized to fit the new, longer cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 --------------------------------------------------------------
# Leading apostrophe keeps the digit-only date code stored as text
# (matching the source report, which never treats 売上日 as a real number).
$ws.Range("B5").Value = "'20210219"
$ws.Range("C5").Value = "高圧化工"
$ws.Range("D5").Value = "ＵＶ－ＰＰＮｏ３ Ａ液 (15KG)"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "CN"
$ws.Range("H5").Value = "1090000840-1"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = "大阪直送"

# --- Row 6 --------------------------------------------------------------
$ws.Range("B6").Value = "'20210219"
$ws.Range("C6").Value = "高圧化工"
$ws.Range("D6").Value = "ＵＶ－ＰＰＮｏ３ Ｂ液 (3.1KG)"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = "CN"
$ws.Range("H6").Value = "1090000840-2"
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = "大阪直送"

# --- Row 7 is no longer shipped out of 本社 on this report; drop it -----
$ws.Rows.Item(7).Delete()

# --- Column widths to accommodate the new text lengths -------------------
# (this host quantizes ColumnWidth to 1/6-character steps, so these are
# the closest achievable values to the target sheet's column widths)
$ws.Columns.Item(3).ColumnWidth = 43 / 6
$ws.Columns.Item(4).ColumnWidth = 170 / 6
$ws.Columns.Item(5).ColumnWidth = 38 / 6
$ws.Columns.Item(8).ColumnWidth = 78 / 6
$ws.Columns.Item(10).ColumnWidth = 43 / 6
